# The document has two floating (anchored) pictures living as two runs
# inside the same paragraph:
#   Shapes.Item(1) -> "Рисунок 1" (docPr id="1", relativeHeight="251661312")
#   Shapes.Item(2) -> "Рисунок 4" (docPr id="4", relativeHeight="251659264")
#
# The edit:
#   1. Nudges the first picture down slightly (positionV posOffset
#      3810 EMU -> 22860 EMU, i.e. 0.3pt -> 1.8pt).
#   2. Removes the second picture entirely.

$d = $word.ActiveDocument

# 1) Move the first picture down (positionV relativeFrom="paragraph").
$pic1 = $d.Shapes.Item(1)
$pic1.Top = 1.8

# 2) Delete the second picture completely.
$pic2 = $d.Shapes.Item(2)
$pic2.Delete()
